$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value2 = 13350.111
$ws.Range("I9").Value2 = 13350.111
$ws.Range("K9").Value2 = 13350.111
$ws.Range("M9").Value2 = -13181.111
$ws.Range("H38").Value2 = 1246.6364
$ws.Range("I38").Value2 = 671.4
$ws.Range("J38").Value2 = 6999
$ws.Range("K38").Value2 = 2014.2
$ws.Range("L38").Value2 = 20997
$ws.Range("M38").Value2 = -1642.2
$ws.Range("N38").Value2 = -21741
$ws.Range("H40").Value2 = 6444.4443
$ws.Range("I40").Value2 = 3571.5715
$ws.Range("J40").Value2 = 8272.637000000001
$ws.Range("K40").Value2 = 3571.5715
$ws.Range("L40").Value2 = 8272.637000000001
$ws.Range("M40").Value2 = -3396.5715
$ws.Range("N40").Value2 = -8622.637000000001
$ws.Range("H43").Value2 = 2358.9412
$ws.Range("I43").Value2 = 2381.4167
$ws.Range("J43").Value2 = 2305
$ws.Range("K43").Value2 = 2381.4167
$ws.Range("L43").Value2 = 2305
$ws.Range("M43").Value2 = -2312.4167
$ws.Range("N43").Value2 = -2443
$ws.Range("H58").Value2 = 8224.083000000001
$ws.Range("J58").Value2 = 9771.223
$ws.Range("L58").Value2 = 29313.669
$ws.Range("N58").Value2 = -29613.669
$ws.Range("H69").Value2 = 17000
$ws.Range("J69").Value2 = 17000
$ws.Range("L69").Value2 = 51000
$ws.Range("N69").Value2 = -52748
$ws.Range("H72").Value2 = 17000
$ws.Range("J72").Value2 = 17000
$ws.Range("L72").Value2 = 153000
$ws.Range("N72").Value2 = -161736
$ws.Range("H76").Value2 = 91014180
$ws.Range("J76").Value2 = 500037500
$ws.Range("L76").Value2 = 500037500
$ws.Range("N76").Value2 = -500038130
$ws.Range("H79").Value2 = 91014180
$ws.Range("J79").Value2 = 500037500
$ws.Range("L79").Value2 = 500037500
$ws.Range("N79").Value2 = -500039684
$ws.Range("H121").Value2 = 5499.5
$ws.Range("J121").Value2 = 5499.5
$ws.Range("L121").Value2 = 16498.5
$ws.Range("N121").Value2 = -19992.5
$ws.Range("H133").Value2 = 58850.785
$ws.Range("J133").Value2 = 58850.785
$ws.Range("L133").Value2 = 58850.785
$ws.Range("N133").Value2 = -68970.785
$ws.Range("H137").Value2 = 4640.4326
$ws.Range("I137").Value2 = 2953.318
$ws.Range("J137").Value2 = 7114.8667
$ws.Range("K137").Value2 = 8859.954000000002
$ws.Range("L137").Value2 = 21344.6001
$ws.Range("M137").Value2 = -6309.954000000002
$ws.Range("N137").Value2 = -26444.6001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value2 = 2430400.8
$ws.Range("I23").Value2 = 2430400.8
$ws.Range("K23").Value2 = 2430400.8
$ws.Range("M23").Value2 = -2430141.8
$ws.Range("H32").Value2 = 2180.3699
$ws.Range("I32").Value2 = 1542.6
$ws.Range("K32").Value2 = 1542.6
$ws.Range("M32").Value2 = -1255.6
$ws.Range("H61").Value2 = 3544.52
$ws.Range("I61").Value2 = 2423.2778
$ws.Range("K61").Value2 = 2423.2778
$ws.Range("M61").Value2 = -2211.2778
$ws.Range("H74").Value2 = 1281.1296
$ws.Range("I74").Value2 = 1082.8298
$ws.Range("K74").Value2 = 1082.8298
$ws.Range("M74").Value2 = -208.8298
$ws.Range("H77").Value2 = 1281.1296
$ws.Range("I77").Value2 = 1082.8298
$ws.Range("K77").Value2 = 5414.148999999999
$ws.Range("M77").Value2 = -1046.148999999999
$ws.Range("H136").Value2 = 3544.52
$ws.Range("I136").Value2 = 2423.2778
$ws.Range("K136").Value2 = 7269.8334
$ws.Range("M136").Value2 = -4719.8334

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value2 = 1199.2858
$ws.Range("I86").Value2 = 1079.6
$ws.Range("K86").Value2 = 1079.6
$ws.Range("M86").Value2 = 43.40000000000009
$ws.Range("H89").Value2 = 1199.2858
$ws.Range("I89").Value2 = 1079.6
$ws.Range("K89").Value2 = 5398
$ws.Range("M89").Value2 = 218

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 2660.3235
$ws.Range("I31").Value2 = 1572.5238
$ws.Range("J31").Value2 = 4417.5386
$ws.Range("K31").Value2 = 1572.5238
$ws.Range("L31").Value2 = 4417.5386
$ws.Range("M31").Value2 = -1277.5238
$ws.Range("N31").Value2 = -5007.5386
$ws.Range("H34").Value2 = 2660.3235
$ws.Range("I34").Value2 = 1572.5238
$ws.Range("J34").Value2 = 4417.5386
$ws.Range("K34").Value2 = 1572.5238
$ws.Range("L34").Value2 = 4417.5386
$ws.Range("M34").Value2 = -1370.5238
$ws.Range("N34").Value2 = -4821.5386
$ws.Range("H132").Value2 = 4237
$ws.Range("I132").Value2 = 2575.6667
$ws.Range("K132").Value2 = 7727.000100000001
$ws.Range("M132").Value2 = -5197.000100000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value2 = 12534
$ws.Range("J39").Value2 = 17148
$ws.Range("L39").Value2 = 51444
$ws.Range("N39").Value2 = -52032
$ws.Range("H131").Value2 = 5024.6113
$ws.Range("I131").Value2 = 879.5
$ws.Range("J131").Value2 = 13314.833
$ws.Range("K131").Value2 = 2638.5
$ws.Range("L131").Value2 = 39944.499
$ws.Range("M131").Value2 = 2401.5
$ws.Range("N131").Value2 = -50024.499
$ws.Range("H138").Value2 = 3000
$ws.Range("I138").Value2 = 1000
$ws.Range("J138").Value2 = 3666.6667
$ws.Range("K138").Value2 = 3000
$ws.Range("L138").Value2 = 11000.0001
$ws.Range("M138").Value2 = 2140
$ws.Range("N138").Value2 = -21280.0001
$ws.Range("H140").Value2 = 4029.6
$ws.Range("I140").Value2 = 3644.111
$ws.Range("K140").Value2 = 10932.333
$ws.Range("M140").Value2 = -5752.332999999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value2 = 2782495
$ws.Range("J3").Value2 = 1434493.6
$ws.Range("L3").Value2 = 1434493.6
$ws.Range("N3").Value2 = -1434725.6
$ws.Range("H10").Value2 = 2014899.8
$ws.Range("I10").Value2 = 3338999.8
$ws.Range("J10").Value2 = 28750
$ws.Range("K10").Value2 = 3338999.8
$ws.Range("L10").Value2 = 28750
$ws.Range("M10").Value2 = -3338830.8
$ws.Range("N10").Value2 = -29088
$ws.Range("H123").Value2 = 43998.5
$ws.Range("J123").Value2 = 43998.5
$ws.Range("L123").Value2 = 43998.5
$ws.Range("H126").Value2 = 90911336
$ws.Range("I126").Value2 = 142858620
$ws.Range("K126").Value2 = 428575860
$ws.Range("M126").Value2 = -428573390
$ws.Range("H132").Value2 = 348786.1
$ws.Range("I132").Value2 = 403731.9
$ws.Range("J132").Value2 = 5374.75
$ws.Range("K132").Value2 = 1211195.7
$ws.Range("L132").Value2 = 16124.25
$ws.Range("M132").Value2 = -1208665.7
$ws.Range("N132").Value2 = -21184.25
$ws.Range("H135").Value2 = 119714.14
$ws.Range("J135").Value2 = 119714.14
$ws.Range("L135").Value2 = 119714.14
$ws.Range("N135").Value2 = -129854.14
$ws.Range("N123").Value2 = -48898.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 596183.3
$ws.Range("I7").Value2 = 917056.0600000001
$ws.Range("K7").Value2 = 917056.0600000001
$ws.Range("M7").Value2 = -916944.0600000001
$ws.Range("H40").Value2 = 1255025.8
$ws.Range("J40").Value2 = 6898.75
$ws.Range("L40").Value2 = 6898.75
$ws.Range("N40").Value2 = -7170.75
$ws.Range("H46").Value2 = 3135.3333
$ws.Range("I46").Value2 = 2452.55
$ws.Range("K46").Value2 = 2452.55
$ws.Range("M46").Value2 = -2264.55
$ws.Range("H61").Value2 = 4856.48
$ws.Range("J61").Value2 = 6908.636
$ws.Range("L61").Value2 = 6908.636
$ws.Range("N61").Value2 = -7312.636
$ws.Range("H68").Value2 = 13737
$ws.Range("I68").Value2 = 14982.667
$ws.Range("J68").Value2 = 10000
$ws.Range("K68").Value2 = 14982.667
$ws.Range("L68").Value2 = 10000
$ws.Range("M68").Value2 = -14233.667
$ws.Range("H71").Value2 = 13737
$ws.Range("I71").Value2 = 14982.667
$ws.Range("J71").Value2 = 10000
$ws.Range("K71").Value2 = 74913.33499999999
$ws.Range("L71").Value2 = 50000
$ws.Range("M71").Value2 = -71169.33499999999
$ws.Range("H113").Value2 = 4856.48
$ws.Range("J113").Value2 = 6908.636
$ws.Range("L113").Value2 = 6908.636
$ws.Range("N113").Value2 = -11248.636
$ws.Range("H126").Value2 = 596183.3
$ws.Range("I126").Value2 = 917056.0600000001
$ws.Range("K126").Value2 = 2751168.18
$ws.Range("M126").Value2 = -2748698.18
$ws.Range("H132").Value2 = 5296.2
$ws.Range("I132").Value2 = 4206.2856
$ws.Range("K132").Value2 = 12618.8568
$ws.Range("M132").Value2 = -10088.8568
$ws.Range("N68").Value2 = -11498
$ws.Range("N71").Value2 = -57488

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value2 = 9940
$ws.Range("J81").Value2 = 16660.223
$ws.Range("L81").Value2 = 33320.446
$ws.Range("N81").Value2 = -35442.446
$ws.Range("H84").Value2 = 9940
$ws.Range("J84").Value2 = 16660.223
$ws.Range("L84").Value2 = 166602.23
$ws.Range("N84").Value2 = -177210.23
$ws.Range("H126").Value2 = 5284.857
$ws.Range("I126").Value2 = 4998.3335
$ws.Range("K126").Value2 = 14995.0005
$ws.Range("M126").Value2 = -12525.0005
